# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Membrillo" (Vega Modelo de Temuco)
# right before the current row 270, pushing the existing rows 270-281 down
# to 271-282 (dimension grows from A1:T281 to A1:T282).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 270..281 down one row, leaving a fresh blank row 270.
$ws.Rows.Item(270).Insert()

# Populate the new row 270 with the new weekly record.
$ws.Range("A270").Value = 10
$ws.Range("B270").Value = "Vega Modelo de Temuco"
$ws.Range("C270").Value = "La Araucanía"
$ws.Range("D270").Value = 45075
$ws.Range("E270").Value = 9
$ws.Range("F270").Value = "Fruta"
$ws.Range("G270").Value = 100104
$ws.Range("H270").Value = "Frutos de pepita"
$ws.Range("I270").Value = 100104003
$ws.Range("J270").Value = "Membrillo"
$ws.Range("K270").Value = "Champion"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 120
$ws.Range("N270").Value = 14000
$ws.Range("O270").Value = 14000
$ws.Range("P270").Value = 14000
$ws.Range("Q270").Value = "$/bandeja 18 kilos granel"
$ws.Range("R270").Value = "Región de O'Higgins"
$ws.Range("S270").Value = 778
$ws.Range("T270").Value = 18
